# SOP_ScrollSaw.docx edit:
# "Prototyping Labs" -> "Prototyping Lab" (singular) in two places:
#   1. The Author/Title/Date table row in the document body footer block.
#   2. The big title in header2 ("Prototyping Labs at GIX" -> "Prototyping Lab at GIX").

$d = $word.ActiveDocument

# 1. Document body: "Title: Prototyping Labs " -> "Title: Prototyping Lab "
$d.Content.Find.Execute("Title: Prototyping Labs ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Title: Prototyping Lab ", 2) | Out-Null

# 2. Headers: "Prototyping Labs at " -> "Prototyping Lab at "
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("Prototyping Labs at ", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "Prototyping Lab at ", 2) | Out-Null
        }
    }
}
